$d = $word.ActiveDocument

# 1) "Целочисленная константа " (spread across two runs) -> "Числовая константа"
$d.Content.Find.Execute("Целочисленная константа ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Числовая константа", 2) | Out-Null

# 2) "Вещественная константа" -> "Строковая константа"
$d.Content.Find.Execute("Вещественная константа", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Строковая константа", 2) | Out-Null

# 3) Move the "_GoBack" bookmark from the trailing empty paragraph to sit right
#    after "Унарная операция" (end of that paragraph's text).
$r = $d.Content
$found = $r.Find.Execute("Унарная операция", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
}
